$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 values (quarter 01-01-2021 revised) ---
$ws.Range("K74").Value = -166
$ws.Range("M74").Value = -97
$ws.Range("P74").Value = -144
$ws.Range("Q74").Value = -144
$ws.Range("S74").Value = -12
$ws.Range("V74").Value = -10
$ws.Range("W74").Value = 1696
$ws.Range("X74").Value = 1040
$ws.Range("Y74").Value = 2
$ws.Range("AB74").Value = 666
$ws.Range("AC74").Value = 1778
$ws.Range("AG74").Value = 66
$ws.Range("AH74").Value = 108
$ws.Range("AI74").Value = 0
$ws.Range("AJ74").Value = -15
$ws.Range("AK74").Value = 123

# --- Append new row 75 (new quarter 01-04-2021) ---
# Force the date label into the sheet as text (matching the existing
# "Serie" column, which stores these labels as plain strings, not as
# real Excel dates), then drop back to the default (unstyled) cell
# format so no visible formatting change is introduced.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"

$ws.Range("B75").Value = 0
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 6428
$ws.Range("G75").Value = 6092
$ws.Range("H75").Value = -107
$ws.Range("I75").Value = 673
$ws.Range("J75").Value = -229
$ws.Range("K75").Value = -58
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 275
$ws.Range("N75").Value = 0
$ws.Range("O75").Value = -333
$ws.Range("P75").Value = 119
$ws.Range("Q75").Value = 119
$ws.Range("R75").Value = 0
$ws.Range("S75").Value = -21
$ws.Range("T75").Value = 4
$ws.Range("U75").Value = -7
$ws.Range("V75").Value = -17
$ws.Range("W75").Value = 303
$ws.Range("X75").Value = 1100
$ws.Range("Y75").Value = 35
$ws.Range("Z75").Value = 76
$ws.Range("AA75").Value = -56
$ws.Range("AB75").Value = -851
$ws.Range("AC75").Value = 2517
$ws.Range("AD75").Value = 1746
$ws.Range("AE75").Value = -47
$ws.Range("AF75").Value = 60
$ws.Range("AG75").Value = 759
$ws.Range("AH75").Value = 120
$ws.Range("AI75").Value = 0
$ws.Range("AJ75").Value = 64
$ws.Range("AK75").Value = 56
